$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (A2:A11, C2:C11, D2:D11) were cyclically shifted down by one:
# the last row of the block (row 11) moved to the top (row 2), pushing the
# rest down by one row. Re-write each changed cell with its new value so the
# table reads correctly (values are pulled from the existing shared-string
# pool, so no new strings are minted).

$ws.Range("C2").Value  = "HEM Aradhana Gangajal"
$ws.Range("D2").Value  = "Aradhana Pure Gangajal 500ml Pack"

$ws.Range("A3").Value  = "Pooja Samagri "
$ws.Range("C3").Value  = "Pooja Samagri Kit "
$ws.Range("D3").Value  = "Aradhana Pooja Samagri"

$ws.Range("C4").Value  = "Aradhana Pooja Oil Pack 900ml"
$ws.Range("D4").Value  = "Aradhana Pooja Oil Pack 20 Pcs"

$ws.Range("A5").Value  = " Pooja Oil "
$ws.Range("C5").Value  = "Aradhana Pooja Oil Pack 450Ml "
$ws.Range("D5").Value  = "Aradhana Pooja Oil Pack 450Ml"

$ws.Range("C6").Value  = "HEM Aradhana Camphor - 100g"
$ws.Range("D6").Value  = "Aradhana 100 Gm Camphor Jar"

$ws.Range("A7").Value  = "HEM  Camphor "
$ws.Range("C7").Value  = "HEM Aradhana Bhimseni Camphor"
$ws.Range("D7").Value  = "Aradhana Bhimseni Camphor (Mrp)"

$ws.Range("C8").Value  = "HEM Aradhana Pure Cow Ghee Diya 100 pcs"
$ws.Range("D8").Value  = "Aradhana Pure Cow Ghee (100 Diyas)"

$ws.Range("C9").Value  = "HEM Aradhana Pure Ghee Diya 100 pcs"
$ws.Range("D9").Value  = "Aradhana Pure Ghee (100 Diyas)"

$ws.Range("C10").Value = "Aradhana Pure Cow Ghee 50 pcs"
$ws.Range("D10").Value = "Aradhana Pure Cow Ghee (50 Diyas)"

$ws.Range("A11").Value = "Ghee Diya"
$ws.Range("C11").Value = "Aradhana Pure Ghee 50 pcs"
$ws.Range("D11").Value = "Aradhana Pure Ghee (50 Diyas)"

# A filter database was defined over the table (Data > Filter), recorded as a
# hidden, sheet-scoped workbook name - without leaving a visible AutoFilter
# button row in the sheet itself.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$14")
$fdb.Visible = $false

# Selection ended up on D18.
$ws.Range("D18").Select()
